$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ36383949",
    "summ36496493",
    "summ36637294",
    "summ36786156",
    "summ36934893",
    "summ37066071",
    "summ37210245",
    "summ37361286",
    "summ37507490",
    "summ37658509",
    "summ37807377",
    "summ37954459",
    "summ38107977",
    "summ38260075",
    "summ38411594",
    "summ38547616",
    "summ38694481",
    "summ38844720",
    "summ38998231",
    "summ39228253",
    "summ39376555",
    "summ39525912",
    "summ39680436",
    "summ39829462",
    "summ39992177",
    "summ40143238",
    "summ40290743",
    "summ40431276",
    "summ40567856",
    "summ40719042",
    "summ40869764",
    "summ41009256",
    "summ41149306",
    "summ41300316",
    "summ41449859",
    "summ41598687",
    "summ41745743",
    "summ41900769",
    "summ42048899",
    "summ42193411",
    "summ42361453",
    "summ42532514",
    "summ42680025",
    "summ42829051",
    "summ42972561",
    "summ43119211",
    "summ43257229",
    "summ43395272",
    "summ43536075",
    "summ43697992"
)

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $ws.Name = $newNames[$i]
}
